$d = $word.ActiveDocument

# --------------------------------------------------------------------
# Goal: the first paragraph currently reads "Update 1, 2" as two runs
# (", 2" is bold, sz=22/szCs=22). We need to append a THIRD run with
# text ", 3" carrying the exact same character formatting (b, bCs,
# sz=22, szCs=22), landing right before the paragraph mark, e.g.:
#
#   <w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/>
#   </w:rPr><w:t>, 3</w:t></w:r>
#
# Directly building that run with Font.Bold/Font.BoldBi/Font.Size/
# Font.SizeBi property assignments on a brand-new run leaves out
# <w:szCs> (the complex-script size does not "stick" on a freshly
# inserted run in this host). Also, editing the text of a run that is
# immediately adjacent to an identically-formatted run (e.g. changing
# "2" to "3" in a duplicate of the ", 2" run placed right after it)
# causes the two runs to be silently coalesced into one on save,
# which would NOT match the target (two separate <w:r> elements).
#
# Workaround: build/edit the new ", 3" run in a scratch paragraph at
# the very end of the document (so it is not adjacent to any
# identically-formatted run and therefore never gets merged), copy
# its *entire* formatting+text via FormattedText (which faithfully
# carries every rPr child, including szCs) from the existing ", 2"
# run, fix up its text there, then copy the corrected FormattedText
# into its final home in paragraph 1, and finally delete the scratch
# paragraph again.
# --------------------------------------------------------------------

$firstPara = $d.Paragraphs.Item(1)

# Locate the existing ", 2" run's range dynamically (fully formatted:
# b, bCs, sz=22, szCs=22) instead of hard-coding character offsets.
$commaTwoRange = $firstPara.Range.Duplicate
$commaTwoRange.Find.Execute(", 2", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null

# 1) Add a scratch paragraph at the very end of the document body to do
#    the text surgery somewhere that is not adjacent to identical
#    formatting (so no accidental run-merging happens).
$docEnd = $d.Content.End
$d.Range($docEnd, $docEnd).InsertParagraphAfter()
$scratchStart = $docEnd

# 2) Copy the fully-formatted ", 2" text into the scratch paragraph.
$scratchTarget = $d.Range($scratchStart, $scratchStart)
$scratchTarget.FormattedText = $commaTwoRange.FormattedText

# 3) Fix up the copied text from ", 2" to ", 3" (safe here: its only
#    neighbor is an empty/differently-formatted run, so nothing merges).
$scratchLastChar = $d.Range($scratchStart + 2, $scratchStart + 3)
$scratchLastChar.Text = "3"
$scratchFinal = $d.Range($scratchStart, $scratchStart + 3)

# 4) Copy the corrected, fully-formatted ", 3" text into place right
#    before paragraph 1's trailing paragraph mark.
$insertAt = $firstPara.Range.End - 1
$finalTarget = $d.Range($insertAt, $insertAt)
$finalTarget.FormattedText = $scratchFinal.FormattedText

# 5) Remove the scratch paragraph (its start shifted by the 3 characters
#    we just inserted into paragraph 1).
$scratchNowStart = $scratchStart + 3
$d.Range($scratchNowStart - 1, $d.Content.End).Delete()
